# Auto-generated Excel COM-interop script to apply the diff changes
# (matches commit: Update gh-pages to output generated at 456a3b4)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value2 = 554
$ws.Range("F6").Value2 = 16
$ws.Range("F7").Value2 = 1904
$ws.Range("F8").Value2 = 5013
$ws.Range("F9").Value2 = 1403
$ws.Range("F11").Value2 = 2988
$ws.Range("F14").Value2 = 1245
$ws.Range("F15").Value2 = 4057
$ws.Range("F16").Value2 = 954
$ws.Range("F27").Value2 = 191
$ws.Range("F28").Value2 = 1041
$ws.Range("F30").Value2 = 93
$ws.Range("F31").Value2 = 108
$ws.Range("F32").Value2 = 172
$ws.Range("F33").Value2 = 1560
$ws.Range("F34").Value2 = 2093
$ws.Range("F37").Value2 = 233
$ws.Range("F38").Value2 = 568
$ws.Range("F39").Value2 = 220
$ws.Range("F42").Value2 = 363
$ws.Range("F43").Value2 = 252
$ws.Range("F45").Value2 = 112

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value2 = 668

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value2 = 668
$ws.Range("F7").Value2 = 554
$ws.Range("F8").Value2 = 16
$ws.Range("F9").Value2 = 1904
$ws.Range("F10").Value2 = 5013
$ws.Range("F11").Value2 = 1403
$ws.Range("F14").Value2 = 2988
$ws.Range("F16").Value2 = 1245
$ws.Range("F17").Value2 = 4057
$ws.Range("F18").Value2 = 954
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value2 = "2024-05-25"
$ws.Range("C32").Value2 = "杭州·原神X星铁X绝区零only"
$ws.Range("D32").Value2 = "望江东路333号 杭州瑞莱克斯大酒店"
$ws.Range("E32").Value2 = "2024.05.25 10:00-05.25 17:00"
$ws.Range("F32").Value2 = 191
$ws.Range("G32").Value2 = 60
$ws.Range("H32").Value2 = "https://show.bilibili.com/platform/detail.html?id=82754"
$ws.Range("I32").Value2 = "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"
$ws.Range("C33").Value2 = "杭州·早鸟5折起·《LALALAND爱乐之城》浪漫主题音乐会"
$ws.Range("D33").Value2 = "武林路77号 浙江省文化馆小剧场（原群艺馆小剧场）"
$ws.Range("E33").Value2 = "2024.05.25 19:30-05.25 21:00"
$ws.Range("F33").Value2 = 4
$ws.Range("G33").Value2 = 100
$ws.Range("H33").Value2 = "https://show.bilibili.com/platform/detail.html?id=84519"
$ws.Range("I33").Value2 = "//i1.hdslb.com/bfs/openplatform/202404/jJLft5tT1712888683239.jpeg"
$ws.Range("F34").Value2 = 1041
$ws.Range("F36").Value2 = 1560
$ws.Range("F37").Value2 = 2093
$ws.Range("F42").Value2 = 233
$ws.Range("F43").Value2 = 568
$ws.Range("F44").Value2 = 220
$ws.Range("F46").Value2 = 363
$ws.Range("F47").Value2 = 252
$ws.Range("F49").Value2 = 112

